$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire column U ("Storage scan Qty" header / sort_qty placeholder).
# This shifts columns V and W left to U and V, matching the target layout.
$ws.Columns("U").Delete()

# Restore the selection reported in the edited file.
$ws.Range("P13").Select()
